$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two vision-statement cells that still need to be refined with the
# text "Revamp this", formatted with Excel's built-in "Bad" cell style
# (dark red text on a light red/pink fill).
$ws.Range("C6").Value = "Revamp this"
$ws.Range("C6").Style = "Bad"

$ws.Range("C12").Value = "Revamp this"
$ws.Range("C12").Style = "Bad"

# Size the new column like the existing text columns.
$ws.Columns.Item(3).ColumnWidth = 11.85546875

# Move the active selection to account for the new column, like the source.
$ws.Range("C13").Select()
